{"js": "// Applies the LOQ4216.docx content updates described by the diff:\n//  1. Ativa\u00e7\u00e3o date bump\n//  2. Docente respons\u00e1vel change\n//  3. M\u00e9todo de avalia\u00e7\u00e3o text\n//  4. Crit\u00e9rio de avalia\u00e7\u00e3o text\n//  5. Norma de recupera\u00e7\u00e3o text\n//  6. Bibliografia text\n//\n// Each change is located via Body.search (exact, case-sensitive match on the\n// original text) and swapped in place with Range.insertText(..., \"Replace\").\n\nasync function replaceText(body, findText, newText) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\nconst body = context.document.body;\n\n// 1. Ativa\u00e7\u00e3o date\nawait replaceText(body, \"Ativa\u00e7\u00e3o: 01/01/2018\", \"Ativa\u00e7\u00e3o: 01/01/2021\");\n\n// 2. Docente(s) Respons\u00e1vel(eis)\nawait replaceText(\n  body,\n  \"5840560 - Marco Antonio Carvalho Pereira\",\n  \"5701460 - Antonio Iacono\"\n);\n\n// 3. M\u00e9todo de avalia\u00e7\u00e3o\nawait replaceText(\n  body,\n  \"Aulas expositivas e pr\u00e1ticas.\",\n  \"Provas e Trabalhos\"\n);\n\n// 4. Crit\u00e9rio de avalia\u00e7\u00e3o\nawait replaceText(\n  body,\n  \"Exerc\u00edcios de aprendizado e exerc\u00edcios de avalia\u00e7\u00e3o far\u00e3o parte da composi\u00e7\u00e3o de notas individuais (NI), com aplica\u00e7\u00e3o de trabalhos pr\u00e1ticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2\",\n  \"M = (0,6P + 0,4T)P = Prova escritaT = Trabalho sobre projeto de f\u00e1bricaM = M\u00e9dia de aproveitamento do alunoAprova\u00e7\u00e3o com m\u00e9dia de aproveitamento maior ou igual a 5,0 e no m\u00ednimo 70% de frequ\u00eancia \u00e0s aulas.A m\u00e9dia das provas deve ser maior ou igual a 5,0 (cinco) para que o aluno possa utilizar a nota do Trabalho.\"\n);\n\n// 5. Norma de recupera\u00e7\u00e3o\nawait replaceText(\n  body,\n  \"A recupera\u00e7\u00e3o dever\u00e1 consistir de uma prova englobando a mat\u00e9ria toda do semestre. - A m\u00e9dia final (p\u00f3s-recupera\u00e7\u00e3o) dever\u00e1 ser composta por uma m\u00e9dia simples entre a nota do semestre (nota final) e a da prova de recupera\u00e7\u00e3o\",\n  \"MF = (0,5 M + 0,5 R)M = M\u00e9dia de aproveitamento do aluno, antes da recupera\u00e7\u00e3oR = Nota de uma prova de recupera\u00e7\u00e3oMF = nota final de aproveitamento, ap\u00f3s a recupera\u00e7\u00e3oAprova\u00e7\u00e3o com m\u00e9dia final de aproveitamento maior ou igual a 5,0.A recupera\u00e7\u00e3o dever\u00e1 consistir de uma prova escrita englobando a mat\u00e9ria toda do semestre.Ter\u00e1 direito \u00e0 prova de recupera\u00e7\u00e3o aqueles alunos reprovados com nota acima de 3,0 e frequ\u00eancia m\u00ednima de 70%.\"\n);\n\n// 6. Bibliografia\nawait replaceText(\n  body,\n  \"M\u00fcther, R. Planejamento do Layout: Sistema SLP. S\u00e3o Paulo, Edgard Bl\u00fccher, 1978. Francischini, P.G.; Gurgel, F.A.C. Administra\u00e7\u00e3o de recursos materiais e patrimoniais. 2a. Edi\u00e7\u00e3o. S\u00e3o Paulo. Editora Cengage. 2013. Slack, N. et al Administra\u00e7\u00e3o da Produ\u00e7\u00e3o. 3\u00aa ed. S\u00e3o Paulo, Atlas, 2009. Valle, C.E. Implanta\u00e7\u00e3o de Ind\u00fastrias. Rio de Janeiro, LTC Editora, 1975.\",\n  \"BANZATO, Eduardo et al. Atualidades na armazenagem. S\u00e3o Paulo: IMAM, 2003.BARNES, R.M. Estudo de Movimentos de Tempos: projeto e medida do trabalho. S\u00e3o Paulo, Edgar Bl\u00fccher, 1977.GURGEL, F.A.C. Administra\u00e7\u00e3o de recursos materiais e patrimoniais. 2a. Edi\u00e7\u00e3o. S\u00e3o Paulo. Editora Cengage. 2013. FRANCISCHINI, P.G.; VALLE, C.E. Implanta\u00e7\u00e3o de Ind\u00fastrias. Rio de Janeiro, LTC Editora, 1975.LEE, Q et al. Projeto de Instala\u00e7\u00f5es e Locais de Trabalho. S\u00e3o Paulo: IMAM, 1998.MOURA, Reinaldo Aparecido. Sistemas e t\u00e9cnicas de movimenta\u00e7\u00e3o e armazenagem de materiais. IMAM, 2012.NEWMANN, C.; SCALICE, R.K. Projeto de F\u00e1brica e Layout. Rio de Janeiro, Elsevier, 2015.M\u00fcther, R. Planejamento do Layout: Sistema SLP. S\u00e3o Paulo, Edgard Bl\u00fccher, 1978. SLACK, Nigel et al. Administra\u00e7\u00e3o da produ\u00e7\u00e3o. S\u00e3o Paulo: Atlas, 8\u00aa ed. 2018.TOMPKINS, James A. et al. Planejamento de instala\u00e7\u00f5es. Editora LTC:, 2013.\"\n);\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# wdFindContinue = 1, wdReplaceAll = 2\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\nfunction Replace-Text($findText, $replaceText) {\n  $range = $d.Content\n  $range.Find.ClearFormatting()\n  $range.Find.Replacement.ClearFormatting()\n  $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $replaceText, $wdReplaceAll)\n}\n\n# 1. Ativa\u00e7\u00e3o date\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2018\" \"Ativa\u00e7\u00e3o: 01/01/2021\"\n\n# 2. Docente respons\u00e1vel\nReplace-Text \"5840560 - Marco Antonio Carvalho Pereira\" \"5701460 - Antonio Iacono\"\n\n# 3. M\u00e9todo de avalia\u00e7\u00e3o\nReplace-Text \"Aulas expositivas e pr\u00e1ticas.\" \"Provas e Trabalhos\"\n\n# 4. Crit\u00e9rio de avalia\u00e7\u00e3o\nReplace-Text \"Exerc\u00edcios de aprendizado e exerc\u00edcios de avalia\u00e7\u00e3o far\u00e3o parte da composi\u00e7\u00e3o de notas individuais (NI), com aplica\u00e7\u00e3o de trabalhos pr\u00e1ticos em grupo (NG). Sendo: Nota Final = (NI+NG)/2\" \"M = (0,6P + 0,4T)P = Prova escritaT = Trabalho sobre projeto de f\u00e1bricaM = M\u00e9dia de aproveitamento do alunoAprova\u00e7\u00e3o com m\u00e9dia de aproveitamento maior ou igual a 5,0 e no m\u00ednimo 70% de frequ\u00eancia \u00e0s aulas.A m\u00e9dia das provas deve ser maior ou igual a 5,0 (cinco) para que o aluno possa utilizar a nota do Trabalho.\"\n\n# 5. Norma de recupera\u00e7\u00e3o\nReplace-Text \"A recupera\u00e7\u00e3o dever\u00e1 consistir de uma prova englobando a mat\u00e9ria toda do semestre. - A m\u00e9dia final (p\u00f3s-recupera\u00e7\u00e3o) dever\u00e1 ser composta por uma m\u00e9dia simples entre a nota do semestre (nota final) e a da prova de recupera\u00e7\u00e3o\" \"MF = (0,5 M + 0,5 R)M = M\u00e9dia de aproveitamento do aluno, antes da recupera\u00e7\u00e3oR = Nota de uma prova de recupera\u00e7\u00e3oMF = nota final de aproveitamento, ap\u00f3s a recupera\u00e7\u00e3oAprova\u00e7\u00e3o com m\u00e9dia final de aproveitamento maior ou igual a 5,0.A recupera\u00e7\u00e3o dever\u00e1 consistir de uma prova escrita englobando a mat\u00e9ria toda do semestre.Ter\u00e1 direito \u00e0 prova de recupera\u00e7\u00e3o aqueles alunos reprovados com nota acima de 3,0 e frequ\u00eancia m\u00ednima de 70%.\"\n\n# 6. Bibliografia\nReplace-Text \"M\u00fcther, R. Planejamento do Layout: Sistema SLP. S\u00e3o Paulo, Edgard Bl\u00fccher, 1978. Francischini, P.G.; Gurgel, F.A.C. Administra\u00e7\u00e3o de recursos materiais e patrimoniais. 2a. Edi\u00e7\u00e3o. S\u00e3o Paulo. Editora Cengage. 2013. Slack, N. et al Administra\u00e7\u00e3o da Produ\u00e7\u00e3o. 3\u00aa ed. S\u00e3o Paulo, Atlas, 2009. Valle, C.E. Implanta\u00e7\u00e3o de Ind\u00fastrias. Rio de Janeiro, LTC Editora, 1975.\" \"BANZATO, Eduardo et al. Atualidades na armazenagem. S\u00e3o Paulo: IMAM, 2003.BARNES, R.M. Estudo de Movimentos de Tempos: projeto e medida do trabalho. S\u00e3o Paulo, Edgar Bl\u00fccher, 1977.GURGEL, F.A.C. Administra\u00e7\u00e3o de recursos materiais e patrimoniais. 2a. Edi\u00e7\u00e3o. S\u00e3o Paulo. Editora Cengage. 2013. FRANCISCHINI, P.G.; VALLE, C.E. Implanta\u00e7\u00e3o de Ind\u00fastrias. Rio de Janeiro, LTC Editora, 1975.LEE, Q et al. Projeto de Instala\u00e7\u00f5es e Locais de Trabalho. S\u00e3o Paulo: IMAM, 1998.MOURA, Reinaldo Aparecido. Sistemas e t\u00e9cnicas de movimenta\u00e7\u00e3o e armazenagem de materiais. IMAM, 2012.NEWMANN, C.; SCALICE, R.K. Projeto de F\u00e1brica e Layout. Rio de Janeiro, Elsevier, 2015.M\u00fcther, R. Planejamento do Layout: Sistema SLP. S\u00e3o Paulo, Edgard Bl\u00fccher, 1978. SLACK, Nigel et al. Administra\u00e7\u00e3o da produ\u00e7\u00e3o. S\u00e3o Paulo: Atlas, 8\u00aa ed. 2018.TOMPKINS, James A. et al. Planejamento de instala\u00e7\u00f5es. Editora LTC:, 2013.\"\n\n$d.Save()\n"}
